$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("A2").Value = 22
$ws.Range("B2").Value = 19

$ws.Range("A3").Value = 11
# B3 unchanged (17)

# A4 unchanged (21)
$ws.Range("B4").Value = 15

# Add new row 5
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 6

# Match the formatting of the other A-column cells (bold/border/center) on the new A5 cell
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
